# Scheduled market-data refresh: update computed price/profit columns (H:N)
# across the per-Leve tables on each job-class sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 10000
$ws.Range("I18").Value = 10000
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 10000
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -9716
$ws.Range("N18").ClearContents()

$ws.Range("H33").Value = 276.46155
$ws.Range("I33").Value = 165.4
$ws.Range("K33").Value = 165.4
$ws.Range("M33").Value = 63.59999999999999

$ws.Range("H40").Value = 2654
$ws.Range("J40").Value = 3249.125
$ws.Range("L40").Value = 3249.125
$ws.Range("N40").Value = -3599.125

$ws.Range("H125").Value = 2375.3333
$ws.Range("I125").Value = 2667.125
$ws.Range("K125").Value = 24004.125
$ws.Range("M125").Value = -21544.125

$ws.Range("H132").Value = 2766.3015
$ws.Range("I132").Value = 2616.7
$ws.Range("K132").Value = 7850.099999999999
$ws.Range("M132").Value = -5320.099999999999

$ws.Range("H135").Value = 2082.889
$ws.Range("I135").Value = 1800.6364
$ws.Range("K135").Value = 16205.7276
$ws.Range("M135").Value = -13670.7276

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4111.8823
$ws.Range("I61").Value = 3892.4614
$ws.Range("K61").Value = 3892.4614
$ws.Range("M61").Value = -3680.4614

$ws.Range("H122").Value = 2905.75
$ws.Range("I122").Value = 2941.3333
$ws.Range("J122").Value = 2799
$ws.Range("K122").Value = 8823.999899999999
$ws.Range("L122").Value = 8397
$ws.Range("M122").Value = -6373.999899999999
$ws.Range("N122").Value = -13297

$ws.Range("H136").Value = 4111.8823
$ws.Range("I136").Value = 3892.4614
$ws.Range("K136").Value = 11677.3842
$ws.Range("M136").Value = -9127.3842

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2616.7273
$ws.Range("I20").Value = 2578.4
$ws.Range("J20").Value = 3000
$ws.Range("K20").Value = 2578.4
$ws.Range("L20").Value = 3000
$ws.Range("M20").Value = -2331.4
$ws.Range("N20").Value = -3494

$ws.Range("H25").Value = 1981
$ws.Range("I25").Value = 1981
$ws.Range("K25").Value = 1981
$ws.Range("M25").Value = -1746

$ws.Range("H86").Value = 7798
$ws.Range("I86").Value = 2871.6875
$ws.Range("K86").Value = 2871.6875
$ws.Range("M86").Value = -1748.6875

$ws.Range("H89").Value = 7798
$ws.Range("I89").Value = 2871.6875
$ws.Range("K89").Value = 14358.4375
$ws.Range("M89").Value = -8742.4375

$ws.Range("H134").Value = 7870.2617
$ws.Range("I134").Value = 7399.727
$ws.Range("K134").Value = 22199.181
$ws.Range("M134").Value = -19664.181

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 627.375
$ws.Range("J22").Value = 787.5
$ws.Range("L22").Value = 787.5
$ws.Range("N22").Value = -1487.5

$ws.Range("H31").Value = 4060.647
$ws.Range("I31").Value = 4128.875
$ws.Range("J31").Value = 4000
$ws.Range("K31").Value = 4128.875
$ws.Range("L31").Value = 4000
$ws.Range("M31").Value = -3833.875
$ws.Range("N31").Value = -4590

$ws.Range("H34").Value = 4060.647
$ws.Range("I34").Value = 4128.875
$ws.Range("J34").Value = 4000
$ws.Range("K34").Value = 4128.875
$ws.Range("L34").Value = 4000
$ws.Range("M34").Value = -3926.875
$ws.Range("N34").Value = -4404

$ws.Range("H41").Value = 49989
$ws.Range("J41").Value = 49989
$ws.Range("L41").Value = 49989
$ws.Range("N41").Value = -50845

$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()

$ws.Range("H59").Value = 90000
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()

$ws.Range("H60").Value = 40046.5
$ws.Range("J60").Value = 40000
$ws.Range("L60").Value = 40000
$ws.Range("N60").Value = -41022

$ws.Range("H62").Value = 24408.2
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 24408.2
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 24408.2
$ws.Range("N62").Value = -25656.2
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 24408.2
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 24408.2
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 122041
$ws.Range("N65").Value = -128281
$ws.Range("M65").ClearContents()

$ws.Range("H68").Value = 40000
$ws.Range("I68").Value = 40000
$ws.Range("K68").Value = 40000
$ws.Range("M68").Value = -39251

$ws.Range("H71").Value = 40000
$ws.Range("I71").Value = 40000
$ws.Range("K71").Value = 120000
$ws.Range("M71").Value = -116256

$ws.Range("H74").Value = 49999
$ws.Range("J74").Value = 49999
$ws.Range("L74").Value = 49999
$ws.Range("N74").Value = -51747

$ws.Range("H77").Value = 49999
$ws.Range("J77").Value = 49999
$ws.Range("L77").Value = 149997
$ws.Range("N77").Value = -158733

$ws.Range("H93").Value = 30796.545
$ws.Range("J93").Value = 28973.75
$ws.Range("L93").Value = 28973.75
$ws.Range("N93").Value = -32717.75

$ws.Range("H122").Value = 3861.1538
$ws.Range("I122").Value = 2489.5557
$ws.Range("K122").Value = 7468.6671
$ws.Range("M122").Value = -5018.6671

$ws.Range("H134").Value = 2379.9092
$ws.Range("I134").Value = 2226.7
$ws.Range("K134").Value = 6680.099999999999
$ws.Range("M134").Value = -4145.099999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1249.2222
$ws.Range("I122").Value = 248.6
$ws.Range("K122").Value = 2237.4
$ws.Range("M122").Value = 212.5999999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8408
$ws.Range("I70").Value = 7700
$ws.Range("J70").Value = 8496.5
$ws.Range("K70").Value = 7700
$ws.Range("L70").Value = 8496.5
$ws.Range("M70").Value = -7430
$ws.Range("N70").Value = -9036.5

$ws.Range("H73").Value = 8408
$ws.Range("I73").Value = 7700
$ws.Range("J73").Value = 8496.5
$ws.Range("K73").Value = 7700
$ws.Range("L73").Value = 8496.5
$ws.Range("M73").Value = -6764
$ws.Range("N73").Value = -10368.5

$ws.Range("H80").Value = 1500
$ws.Range("I80").Value = 1500
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 1500
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -502
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 1500
$ws.Range("I83").Value = 1500
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 7500
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -2508
$ws.Range("N83").ClearContents()

$ws.Range("H97").Value = 2630.5625
$ws.Range("I97").Value = 554.63635
$ws.Range("K97").Value = 554.63635
$ws.Range("M97").Value = -58.63634999999999

$ws.Range("H113").Value = 6265.8
$ws.Range("I113").Value = 5938.3335
$ws.Range("K113").Value = 5938.3335
$ws.Range("M113").Value = -3768.3335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 3173.3333
$ws.Range("I9").Value = 3173.3333
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 3173.3333
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -2949.3333
$ws.Range("N9").ClearContents()

$ws.Range("H68").Value = 2478.1
$ws.Range("I68").Value = 2483.2856
$ws.Range("J68").Value = 2466
$ws.Range("K68").Value = 2483.2856
$ws.Range("L68").Value = 2466
$ws.Range("M68").Value = -1734.2856
$ws.Range("N68").Value = -3964

$ws.Range("H71").Value = 2478.1
$ws.Range("I71").Value = 2483.2856
$ws.Range("J71").Value = 2466
$ws.Range("K71").Value = 12416.428
$ws.Range("L71").Value = 12330
$ws.Range("M71").Value = -8672.428
$ws.Range("N71").Value = -19818

$ws.Range("H93").Value = 1239.3334
$ws.Range("I93").Value = 1022.2857
$ws.Range("J93").Value = 1999
$ws.Range("K93").Value = 1022.2857
$ws.Range("L93").Value = 1999
$ws.Range("M93").Value = 225.7143
$ws.Range("N93").Value = -4495

$ws.Range("H122").Value = 5412.5
$ws.Range("I122").Value = 4865.5713
$ws.Range("K122").Value = 14596.7139
$ws.Range("M122").Value = -12146.7139

$ws.Range("H136").Value = 2194.6
$ws.Range("I136").Value = 2026.4166
$ws.Range("K136").Value = 6079.2498
$ws.Range("M136").Value = -3529.2498

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H117").Value = 44998
$ws.Range("J117").Value = 44998
$ws.Range("L117").Value = 44998
$ws.Range("N117").Value = -54176

$ws.Range("H122").Value = 9070.579
$ws.Range("I122").Value = 2762.2666
$ws.Range("J122").Value = 32726.75
$ws.Range("K122").Value = 8286.799800000001
$ws.Range("L122").Value = 98180.25
$ws.Range("M122").Value = -5836.799800000001
$ws.Range("N122").Value = -103080.25

$ws.Range("H126").Value = 2610.4443
$ws.Range("J126").Value = 2799.5
$ws.Range("L126").Value = 8398.5
$ws.Range("N126").Value = -13338.5
